$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 797645.7
$ws.Range("I17").Value = 257.14285
$ws.Range("J17").Value = 916405.7
$ws.Range("K17").Value = 771.4285500000001
$ws.Range("L17").Value = 2749217.1
$ws.Range("M17").Value = -603.4285500000001
$ws.Range("N17").Value = -2749553.1

$ws.Range("H106").Value = 4276505
$ws.Range("I106").Value = 4833845
$ws.Range("K106").Value = 4833845
$ws.Range("M106").Value = -4833214

$ws.Range("H117").Value = 30000
$ws.Range("J117").Value = 30000
$ws.Range("L117").Value = 30000
$ws.Range("N117").Value = -39178

$ws.Range("H132").Value = 252885.58
$ws.Range("I132").Value = 265019.28
$ws.Range("J132").Value = 66835.336
$ws.Range("K132").Value = 795057.8400000001
$ws.Range("L132").Value = 200506.008
$ws.Range("M132").Value = -792527.8400000001
$ws.Range("N132").Value = -205566.008

$ws.Range("H137").Value = 30304386
$ws.Range("I137").Value = 55556548
$ws.Range("J137").Value = 1792.2667
$ws.Range("K137").Value = 166669644
$ws.Range("L137").Value = 5376.800099999999
$ws.Range("M137").Value = -166667094
$ws.Range("N137").Value = -10476.8001

$ws.Range("H138").Value = 4263653.5
$ws.Range("I138").Value = 948476.2
$ws.Range("J138").Value = 7578830.5
$ws.Range("K138").Value = 2845428.6
$ws.Range("L138").Value = 22736491.5
$ws.Range("M138").Value = -2840288.6
$ws.Range("N138").Value = -22746771.5

$ws.Range("H141").Value = 2069.9832
$ws.Range("I141").Value = 1367.5193
$ws.Range("J141").Value = 7288.2856
$ws.Range("K141").Value = 4102.5579
$ws.Range("L141").Value = 21864.8568
$ws.Range("M141").Value = 1077.4421
$ws.Range("N141").Value = -32224.8568


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19219.58
$ws.Range("I32").Value = 2889.2856
$ws.Range("J32").Value = 171635.67
$ws.Range("K32").Value = 2889.2856
$ws.Range("L32").Value = 171635.67
$ws.Range("M32").Value = -2602.2856
$ws.Range("N32").Value = -172209.67

$ws.Range("H45").Value = 1031.5
$ws.Range("I45").Value = 1050.1666
$ws.Range("K45").Value = 1050.1666
$ws.Range("M45").Value = -673.1666

$ws.Range("H74").Value = 4114.8086
$ws.Range("I74").Value = 1250.129
$ws.Range("J74").Value = 9665.125
$ws.Range("K74").Value = 1250.129
$ws.Range("L74").Value = 9665.125
$ws.Range("M74").Value = -376.1289999999999
$ws.Range("N74").Value = -11413.125

$ws.Range("H77").Value = 4114.8086
$ws.Range("I77").Value = 1250.129
$ws.Range("J77").Value = 9665.125
$ws.Range("K77").Value = 6250.645
$ws.Range("L77").Value = 48325.625
$ws.Range("M77").Value = -1882.645
$ws.Range("N77").Value = -57061.625

$ws.Range("H110").Value = 757.1111
$ws.Range("I110").Value = 647.6667
$ws.Range("J110").Value = 1304.3334
$ws.Range("K110").Value = 647.6667
$ws.Range("L110").Value = 1304.3334
$ws.Range("M110").Value = 1397.3333
$ws.Range("N110").Value = -5394.3334

$ws.Range("H122").Value = 2285.12
$ws.Range("I122").Value = 2383.7646
$ws.Range("J122").Value = 2075.5
$ws.Range("K122").Value = 7151.293799999999
$ws.Range("L122").Value = 6226.5
$ws.Range("M122").Value = -4701.293799999999
$ws.Range("N122").Value = -11126.5


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12623.818
$ws.Range("I86").Value = 3981.875
$ws.Range("J86").Value = 35669
$ws.Range("K86").Value = 3981.875
$ws.Range("L86").Value = 35669
$ws.Range("M86").Value = -2858.875
$ws.Range("N86").Value = -37915

$ws.Range("H89").Value = 12623.818
$ws.Range("I89").Value = 3981.875
$ws.Range("J89").Value = 35669
$ws.Range("K89").Value = 19909.375
$ws.Range("L89").Value = 178345
$ws.Range("M89").Value = -14293.375
$ws.Range("N89").Value = -189577

$ws.Range("H105").Value = 3379
$ws.Range("I105").Value = 3281.7058
$ws.Range("J105").Value = 3516.8333
$ws.Range("K105").Value = 3281.7058
$ws.Range("L105").Value = 3516.8333
$ws.Range("M105").Value = -1534.7058
$ws.Range("N105").Value = -7010.8333


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1743.4166
$ws.Range("I31").Value = 1024.7391
$ws.Range("J31").Value = 3014.923
$ws.Range("K31").Value = 1024.7391
$ws.Range("L31").Value = 3014.923
$ws.Range("M31").Value = -729.7391
$ws.Range("N31").Value = -3604.923

$ws.Range("H34").Value = 1743.4166
$ws.Range("I34").Value = 1024.7391
$ws.Range("J34").Value = 3014.923
$ws.Range("K34").Value = 1024.7391
$ws.Range("L34").Value = 3014.923
$ws.Range("M34").Value = -822.7391
$ws.Range("N34").Value = -3418.923

$ws.Range("H99").Value = 7813503
$ws.Range("I99").Value = 12500822
$ws.Range("K99").Value = 12500822
$ws.Range("M99").Value = -12499324

$ws.Range("H105").Value = 913.3125
$ws.Range("I105").Value = 840.93335
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 840.93335
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = 906.06665
$ws.Range("N105").Value = -5493

$ws.Range("H126").Value = 7813503
$ws.Range("I126").Value = 12500822
$ws.Range("K126").Value = 37502466
$ws.Range("M126").Value = -37499996

$ws.Range("H132").Value = 1881.3208
$ws.Range("I132").Value = 1557
$ws.Range("J132").Value = 2879.2307
$ws.Range("K132").Value = 4671
$ws.Range("L132").Value = 8637.6921
$ws.Range("M132").Value = -2141
$ws.Range("N132").Value = -13697.6921

$ws.Range("H134").Value = 2114.1296
$ws.Range("I134").Value = 1479.1025
$ws.Range("K134").Value = 4437.3075
$ws.Range("M134").Value = -1902.3075


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1060.75
$ws.Range("I5").Value = 681.5454999999999
$ws.Range("J5").Value = 1895
$ws.Range("K5").Value = 2044.6365
$ws.Range("L5").Value = 5685
$ws.Range("M5").Value = -1932.6365
$ws.Range("N5").Value = -5909

$ws.Range("H109").Value = 2076
$ws.Range("I109").Value = 1141.1428
$ws.Range("J109").Value = 3166.6667
$ws.Range("K109").Value = 3423.4284
$ws.Range("L109").Value = 9500.000100000001
$ws.Range("M109").Value = -2383.4284
$ws.Range("N109").Value = -11580.0001

$ws.Range("H122").Value = 466.86957
$ws.Range("I122").Value = 261.35294
$ws.Range("J122").Value = 1049.1666
$ws.Range("K122").Value = 2352.17646
$ws.Range("L122").Value = 9442.499400000001
$ws.Range("M122").Value = 97.82354000000032
$ws.Range("N122").Value = -14342.4994

$ws.Range("H124").Value = 1800
$ws.Range("J124").Value = 1800
$ws.Range("L124").Value = 5400
$ws.Range("N124").Value = -15220

$ws.Range("H125").Value = 2860.4062
$ws.Range("J125").Value = 3001.1785
$ws.Range("L125").Value = 9003.5355
$ws.Range("N125").Value = -18843.5355

$ws.Range("H131").Value = 1099.375
$ws.Range("J131").Value = 2400
$ws.Range("L131").Value = 7200
$ws.Range("N131").Value = -17280

$ws.Range("H132").Value = 1511.3334
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 1633.6
$ws.Range("K132").Value = 8100
$ws.Range("L132").Value = 14702.4
$ws.Range("M132").Value = -5570
$ws.Range("N132").Value = -19762.4

$ws.Range("H133").Value = 25799.092
$ws.Range("I133").Value = 2827.1428
$ws.Range("K133").Value = 8481.428400000001
$ws.Range("M133").Value = -3421.428400000001

$ws.Range("H135").Value = 1060.75
$ws.Range("I135").Value = 681.5454999999999
$ws.Range("J135").Value = 1895
$ws.Range("K135").Value = 6133.9095
$ws.Range("L135").Value = 17055
$ws.Range("M135").Value = -3598.9095
$ws.Range("N135").Value = -22125

$ws.Range("H137").Value = 5054366
$ws.Range("I137").Value = 7696239.5
$ws.Range("J137").Value = 148030.42
$ws.Range("K137").Value = 23088718.5
$ws.Range("L137").Value = 444091.26
$ws.Range("M137").Value = -23083618.5
$ws.Range("N137").Value = -454291.26


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2280.9524
$ws.Range("I126").Value = 1700
$ws.Range("K126").Value = 5100
$ws.Range("M126").Value = -2630


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1343.8948
$ws.Range("I16").Value = 1470.75
$ws.Range("J16").Value = 667.3333
$ws.Range("K16").Value = 1470.75
$ws.Range("L16").Value = 667.3333
$ws.Range("M16").Value = -1300.75
$ws.Range("N16").Value = -1007.3333

$ws.Range("H43").Value = 15400
$ws.Range("J43").Value = 6750
$ws.Range("L43").Value = 6750
$ws.Range("N43").Value = -7136

$ws.Range("H93").Value = 1055.5555
$ws.Range("I93").Value = 884.3333
$ws.Range("J93").Value = 1398
$ws.Range("K93").Value = 884.3333
$ws.Range("L93").Value = 1398
$ws.Range("M93").Value = 363.6667
$ws.Range("N93").Value = -3894

$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 11547
$ws.Range("J58").Value = 11547
$ws.Range("L58").Value = 11547
$ws.Range("N58").Value = -12163

$ws.Range("H122").Value = 68760.266
$ws.Range("I122").Value = 85200.336
$ws.Range("K122").Value = 255601.008
$ws.Range("M122").Value = -253151.008

$ws.Range("H132").Value = 7938084.5
$ws.Range("I132").Value = 11629482
$ws.Range("J132").Value = 1579.7
$ws.Range("K132").Value = 34888446
$ws.Range("L132").Value = 4739.1
$ws.Range("M132").Value = -34885916
$ws.Range("N132").Value = -9799.1

$ws.Range("H136").Value = 22320.086
$ws.Range("I136").Value = 23137.31
$ws.Range("J136").Value = 3932.5
$ws.Range("K136").Value = 69411.93000000001
$ws.Range("L136").Value = 11797.5
$ws.Range("M136").Value = -66861.93000000001
$ws.Range("N136").Value = -16897.5

